# ajout de la maitrise
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("donnees")

# Rename the "hero base stat" labels (rows 4-7) to the "base" wording
$ws.Range("A4").Value = "atk base  hero"
$ws.Range("A5").Value = "def base hero"
$ws.Range("A6").Value = "pv base hero"
$ws.Range("A7").Value = "vit base hero"

# Row 8 used to hold "% augmentation stats par level" / "10.0" (text).
# It now becomes "pp base attaque  1" with a numeric value, so the old
# attack-power rows (9-12) shift up by one and get the "base" wording.
$ws.Range("A8").Value = "pp base attaque  1"
$ws.Range("B8").Value = 5

$ws.Range("A9").Value = "pp base attaque  2"
$ws.Range("B9").Value = 6

$ws.Range("A10").Value = "pp base attaque  3"
$ws.Range("B10").Value = 7

$ws.Range("A11").Value = "pp base attaque  4"
$ws.Range("B11").Value = 8

# New "maitrise" (mastery) rows
$ws.Range("A12").Value = "point maitrise gagné par niv"
$ws.Range("B12").Value = 10

$ws.Range("A13").Value = " + atk par pt maitrise"
$ws.Range("B13").Value = 1

$ws.Range("A14").Value = " + def par pt maitrise"
$ws.Range("B14").Value = 5

$ws.Range("A15").Value = " + vit par pt maitrise"
$ws.Range("B15").Value = 2

$ws.Range("A16").Value = " + pv par pt maitrise"
$ws.Range("B16").Value = 1

$ws.Range("A17").Value = " + pp par pt maitrise"
$ws.Range("B17").Value = 1

[void]$ws.Range("B16").Select()
